
$d = $word.ActiveDocument

# 1. Replace "Andy" with "Ryan" (mail merge preview greeting name)
$d.Content.Find.Execute("Andy", $true, $false, $false, $false, $false, $true, 1, $false, "Ryan", 2) | Out-Null

# 2. Insert a new paragraph after the greeting paragraph with the explanatory italic text
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$newParaXml = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t xml:space=`"preserve`">Below is a copy of the type of message along with links to the November Heritage Happenings </w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t>n</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t xml:space=`"preserve`">ewsletter and </w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t>c</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t>alenda</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t xml:space=`"preserve`">r — published by residents of Heritage on the Marina. </w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t>I plan to add your name as a monthly subscriber. At a</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t>n</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t>y time</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t xml:space=`"preserve`"> just shout and your name will be removed</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=`"Leelawadee UI`"/><w:i/><w:iCs/></w:rPr><w:t xml:space=`"preserve`"> from the list.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$r2.InsertXML($newParaXml)

# 3. Split "Novembe" into "Novem" + "be" runs (keep trailing "r" run as-is) in the Newsletter hyperlink
$rngNov = $d.Content
$rngNov.Find.Execute("Novembe") | Out-Null
$splitNov = $d.Range($rngNov.End - 2, $rngNov.End)
$splitNov.Font.Bold = $true

# 4. Split " November " into " November" + " " runs in the Calendar hyperlink (second occurrence)
$start = 0
$foundCount = 0
$targetRng = $null
while ($true) {
    $cur = $d.Range($start, $d.Content.End)
    $found = $cur.Find.Execute(" November ")
    if (-not $found) { break }
    $foundCount = $foundCount + 1
    if ($foundCount -eq 2) {
        $targetRng = $d.Range($cur.Start, $cur.End)
        break
    }
    $start = $cur.End
}
$splitNovCal = $d.Range($targetRng.End - 1, $targetRng.End)
$splitNovCal.Font.Bold = $true

Write-Host "All edits applied"
